# Apply the updated cryptocurrency price/volume snapshot to the sheet.
# Text-like numeric values (e.g. "0.999", "309.81") are written with a leading
# apostrophe so Excel keeps them as text instead of coercing them to numbers,
# matching the original inline-string cell type. The Style reset afterwards
# clears the auto-applied 'quote prefix' formatting so the cell's style index
# is left unchanged (same as before the edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''43.313.49'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.80%  '
$ws.Range("E2").Style = "Normal"
# Row 3
$ws.Range("D3").Value = '''2.356.34'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +3.19%  '
$ws.Range("E3").Style = "Normal"
# Row 4
$ws.Range("D4").Value = '''0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("E4").Style = "Normal"
# Row 5
$ws.Range("D5").Value = '''309.81'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.09%  '
$ws.Range("E5").Style = "Normal"
# Row 6
$ws.Range("D6").Value = '''104.04'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.09%  '
$ws.Range("E6").Style = "Normal"
# Row 7
$ws.Range("E7").Value = '  -0.69%  '
$ws.Range("E7").Style = "Normal"
# Row 8
$ws.Range("D8").Value = '''0.999'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("E8").Style = "Normal"
# Row 9
$ws.Range("E9").Value = '  +3.62%  '
$ws.Range("E9").Style = "Normal"
# Row 10
$ws.Range("D10").Value = '''36.03'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.76%  '
$ws.Range("E10").Style = "Normal"
# Row 11
$ws.Range("D11").Value = '''52.83'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.52%  '
$ws.Range("E11").Style = "Normal"
# Row 12
$ws.Range("E12").Value = '  -1.05%  '
$ws.Range("E12").Style = "Normal"
# Row 13
$ws.Range("E13").Value = '  -0.53%  '
$ws.Range("E13").Style = "Normal"
# Row 14
$ws.Range("D14").Value = '''6.99'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.07%  '
$ws.Range("E14").Style = "Normal"
# Row 15
$ws.Range("D15").Value = '''2.722.29'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.37%  '
$ws.Range("E15").Style = "Normal"
# Row 16
$ws.Range("D16").Value = '''15.75'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +6.67%  '
$ws.Range("E16").Style = "Normal"
# Row 17
$ws.Range("D17").Value = '''2.338.49'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.65%  '
$ws.Range("E17").Style = "Normal"
# Row 18
$ws.Range("D18").Value = '''0.810'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.49%  '
$ws.Range("E18").Style = "Normal"
# Row 19
$ws.Range("D19").Value = '''43.296.12'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.20%  '
$ws.Range("E19").Style = "Normal"
# Row 20
$ws.Range("D20").Value = '''11.94'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.24%  '
$ws.Range("E20").Style = "Normal"
# Row 21
$ws.Range("E21").Value = '  +1.44%  '
$ws.Range("E21").Style = "Normal"
# Row 22
$ws.Range("D22").Value = '''6.27'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.81%  '
$ws.Range("E22").Style = "Normal"
# Row 23
$ws.Range("D23").Value = '''68.25'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.83%  '
$ws.Range("E23").Style = "Normal"
# Row 24
$ws.Range("D24").Value = '''241.57'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.31%  '
$ws.Range("E24").Style = "Normal"
# Row 25
$ws.Range("D25").Value = '''2.05'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.25%  '
$ws.Range("E25").Style = "Normal"
# Row 26
$ws.Range("E26").Value = '  +0.47%  '
$ws.Range("E26").Style = "Normal"
# Row 27
$ws.Range("E27").Value = '  +0.60%  '
$ws.Range("E27").Style = "Normal"
# Row 28
$ws.Range("D28").Value = '''25.62'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +7.81%  '
$ws.Range("E28").Style = "Normal"
# Row 29
$ws.Range("B29").Value = '''InjectiveProtocol'
$ws.Range("B29").Style = "Normal"
$ws.Range("C29").Value = '''https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("C29").Style = "Normal"
$ws.Range("D29").Value = '''36.54'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -5.54%  '
$ws.Range("E29").Style = "Normal"
# Row 30
$ws.Range("B30").Value = '''Toncoin'
$ws.Range("B30").Style = "Normal"
$ws.Range("C30").Value = '''https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("C30").Style = "Normal"
$ws.Range("D30").Value = '''2.21'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.61%  '
$ws.Range("E30").Style = "Normal"
# Row 31
$ws.Range("E31").Value = '  -0.39%  '
$ws.Range("E31").Style = "Normal"
# Row 32
$ws.Range("D32").Value = '''162.67'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.23%  '
$ws.Range("E32").Style = "Normal"
# Row 33
$ws.Range("E33").Value = '  -0.28%  '
$ws.Range("E33").Style = "Normal"
# Row 34
$ws.Range("D34").Value = '''0.999'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.11%  '
$ws.Range("E34").Style = "Normal"
# Row 35
$ws.Range("D35").Value = '''18.26'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.66%  '
$ws.Range("E35").Style = "Normal"
# Row 36
$ws.Range("E36").Value = '  +6.63%  '
$ws.Range("E36").Style = "Normal"
# Row 37
$ws.Range("D37").Value = '''3.13'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.19%  '
$ws.Range("E37").Style = "Normal"
# Row 38
$ws.Range("D38").Value = '''0.0739'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.74%  '
$ws.Range("E38").Style = "Normal"
# Row 39
$ws.Range("B39").Value = '''ARBITRUM'
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").Value = '''https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").Value = '''1.92'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.12%  '
$ws.Range("E39").Style = "Normal"
# Row 40
$ws.Range("B40").Value = '''RenderToken'
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = '''https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = '''4.58'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +8.79%  '
$ws.Range("E40").Style = "Normal"
# Row 41
$ws.Range("D41").Value = '''0.107'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.31%  '
$ws.Range("E41").Style = "Normal"
# Row 42
$ws.Range("E42").Value = '  -0.37%  '
$ws.Range("E42").Style = "Normal"
# Row 43
$ws.Range("D43").Value = '''2.42'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +6.02%  '
$ws.Range("E43").Style = "Normal"
# Row 44
$ws.Range("D44").Value = '''20.04'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.22%  '
$ws.Range("E44").Style = "Normal"
# Row 45
$ws.Range("D45").Value = '''0.0292'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.73%  '
$ws.Range("E45").Style = "Normal"
# Row 46
$ws.Range("D46").Value = '''1.988.17'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.07%  '
$ws.Range("E46").Style = "Normal"
# Row 47
$ws.Range("B47").Value = '''NEARProtocol'
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = '''https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = '''3.08'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.51%  '
$ws.Range("E47").Style = "Normal"
# Row 48
$ws.Range("B48").Value = '''FraxShare'
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = '''https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = '''10.42'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +6.42%  '
$ws.Range("E48").Style = "Normal"
# Row 49
$ws.Range("D49").Value = '''58.77'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +7.50%  '
$ws.Range("E49").Style = "Normal"
# Row 50
$ws.Range("E50").Value = '  +4.85%  '
$ws.Range("E50").Style = "Normal"
# Row 51
$ws.Range("D51").Value = '''2.91'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.55%  '
$ws.Range("E51").Style = "Normal"
